$wb = $excel.ActiveWorkbook

$wsInventory = $wb.Worksheets.Item("Inventory")
$wsProducts  = $wb.Worksheets.Item("Products")

# --- Products sheet: replace the hyperlinked image-URL column with plain
# image-filename text (no hyperlinks), using soft-assert friendly values. ---
$imgRange = $wsProducts.Range("B2:B7")
$imgRange.Hyperlinks.Delete()
$imgRange.Style = "Normal"

$wsProducts.Range("B2").Value = "sauce-backpack-1200x1500"
$wsProducts.Range("B3").Value = "bike-light-1200x1500"
$wsProducts.Range("B4").Value = "bolt-shirt-1200x1500"
$wsProducts.Range("B5").Value = "sauce-pullover-1200x1500"
$wsProducts.Range("B6").Value = "red-onesie-1200x1500"
$wsProducts.Range("B7").Value = "red-tatt-1200x1500"

$imgRange.Interior.ColorIndex = -4142

# --- Selections: move off the Inventory sheet, onto the Products sheet. ---
$wsInventory.Range("B10").Select()

$wsProducts.Activate()
$wsProducts.Range("B8").Select()
$excel.ActiveWindow.Zoom = 133

Write-Host "done"
